# Applies the odds-update edit described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 6
$ws.Range("K3").Value = 1.95
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("S3").Value = 1.57
$ws.Range("T3").Value = 2.25
$ws.Range("AH3").Value = 10
$ws.Range("AI3").Value = 26
$ws.Range("AJ3").Value = 21
$ws.Range("AN3").Value = 3.4
$ws.Range("AQ3").Value = 34
$ws.Range("AT3").Value = 2.25
$ws.Range("AW3").Value = 7

# Row 4 updates
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53

# Row 9 updates
$ws.Range("G9").Value = 1.42
$ws.Range("AC9").Value = 19
$ws.Range("AD9").Value = 9.5
$ws.Range("AH9").Value = 23
$ws.Range("AJ9").Value = 21
$ws.Range("BC9").Value = 151
$ws.Range("BD9").Value = 176

# Row 10 updates
$ws.Range("BD10").Value = 151
